# Update the "Periodo Mora" (column E) values and the corresponding
# "Valor Mora" (column F) values for rows 16-22 on the single worksheet.
# This reflects removing the oldest period (2412) from the top of the
# list and appending new periods (2501-2505) above it, shifting the
# existing rows, while keeping the "Valor Mora" amounts tied to their
# original periods.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$periods = @("2505", "2504", "2503", "2502", "2501", "2412", "2506")
$valores  = @(52000, 52000, 52000, 52000, 52000, 45066, 45066)

for ($i = 0; $i -lt $periods.Length; $i++) {
    $row = 16 + $i
    $ws.Cells.Item($row, 5).Value = $periods[$i]
    $ws.Cells.Item($row, 6).Value = $valores[$i]
}
